$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  "18+7=25",
  "16+65=81",
  "76-44=32",
  "59-21=38",
  "53+13=66",
  "30-14=16",
  "48+46=94",
  "6+18=24",
  "3+55=58",
  "17+19=36",
  "42+53=95",
  "57+33=90",
  "70+11=81",
  "72+14=86",
  "88-9=79",
  "65-33=32",
  "63+32=95",
  "0+88=88",
  "28+29=57",
  "71-58=13",
  "75-41=34",
  "39+8=47",
  "69-17=52",
  "99-82=17",
  "35+33=68",
  "30+1=31",
  "82-64=18",
  "80-46=34",
  "53-37=16",
  "46+29=75",
  "41+9=50",
  "52+0=52",
  "23-18=5",
  "8+70=78",
  "14+28=42",
  "86-4=82",
  "71-15=56",
  "46-3=43",
  "87-57=30",
  "3+83=86",
  "44+43=87",
  "19+20=39",
  "87-18=69",
  "80-19=61",
  "60-40=20",
  "85-28=57",
  "22+21=43",
  "41+23=64",
  "19+27=46",
  "65-29=36",
  "91-56=35",
  "50+21=71",
  "50+17=67",
  "18+24=42",
  "31+53=84",
  "17+61=78",
  "69+3=72",
  "20+3=23",
  "34+58=92",
  "13+76=89",
  "1+73=74",
  "53-7=46",
  "63+26=89",
  "74+2=76",
  "89+5=94",
  "25-14=11",
  "62+7=69",
  "72-1=71",
  "54-30=24",
  "62-32=30",
  "51+45=96",
  "8+85=93",
  "3+90=93",
  "15+7=22",
  "30-17=13",
  "72-65=7",
  "17-11=6",
  "43+26=69",
  "14+1=15",
  "19-16=3",
  "10+61=71",
  "88-82=6",
  "86-0=86",
  "62+27=89",
  "15+50=65",
  "59-35=24",
  "88-44=44",
  "50-41=9",
  "25-0=25",
  "68+16=84",
  "53-32=21",
  "44+21=65",
  "88-72=16",
  "54+18=72",
  "20+19=39",
  "12+21=33",
  "80+7=87",
  "61+12=73",
  "93-18=75",
  "59-26=33"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$idx]
    $idx++
  }
}

Write-Host "Updated" $idx "cells"